# Applies the per-coin Price (D) and Volume(1h) (E) refresh, and the
# TheGraph/Arweave row swap (rows 44-45), per the commit's data update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text so Excel does not reinterpret
# number-like strings (e.g. "1.00", "0.996") as numeric values, which
# would silently drop significant trailing/leading zeros.
function Set-TextValue($cell, $text) {
    $range = $ws.Range($cell)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-TextValue "D2" "68.071.39"
Set-TextValue "E2" "  +0.56%  "
Set-TextValue "D3" "3.789.81"
Set-TextValue "E3" "  -0.25%  "
Set-TextValue "D4" "0.996"
Set-TextValue "E4" "  -0.48%  "
Set-TextValue "D5" "600.02"
Set-TextValue "E5" "  +0.63%  "
Set-TextValue "D6" "164.71"
Set-TextValue "E6" "  -1.45%  "
Set-TextValue "E7" "  +0.13%  "
Set-TextValue "E8" "  -0.57%  "
Set-TextValue "E9" "  -1.26%  "
Set-TextValue "E10" "  +0.00%  "
Set-TextValue "D11" "6.58"
Set-TextValue "E11" "  +4.39%  "
Set-TextValue "D12" "0.0000247"
Set-TextValue "E12" "  -2.28%  "
Set-TextValue "D13" "35.62"
Set-TextValue "E13" "  -0.81%  "
Set-TextValue "D14" "4.428.27"
Set-TextValue "E14" "  -0.27%  "
Set-TextValue "D15" "3.796.02"
Set-TextValue "E15" "  -0.24%  "
Set-TextValue "D16" "68.076.27"
Set-TextValue "E16" "  +0.61%  "
Set-TextValue "D17" "18.26"
Set-TextValue "E17" "  -1.67%  "
Set-TextValue "E18" "  +2.28%  "
Set-TextValue "E19" "  -0.55%  "
Set-TextValue "D20" "460.17"
Set-TextValue "E20" "  -0.18%  "
Set-TextValue "D21" "9.64"
Set-TextValue "E21" "  -2.81%  "
Set-TextValue "E22" "  -0.23%  "
Set-TextValue "E23" "  -4.69%  "
Set-TextValue "D24" "82.97"
Set-TextValue "E24" "  -0.60%  "
Set-TextValue "D25" "11.96"
Set-TextValue "E25" "  -1.11%  "
Set-TextValue "E26" "  -0.34%  "
Set-TextValue "E27" "  -0.11%  "
Set-TextValue "D28" "9.96"
Set-TextValue "E28" "  -0.37%  "
Set-TextValue "D29" "3.940.69"
Set-TextValue "E29" "  +0.02%  "
Set-TextValue "D30" "2.23"
Set-TextValue "E31" "  -5.18%  "
Set-TextValue "D32" "7.28"
Set-TextValue "E32" "  +0.04%  "
Set-TextValue "D33" "29.22"
Set-TextValue "E33" "  -1.13%  "
Set-TextValue "D34" "1.00"
Set-TextValue "D35" "8.98"
Set-TextValue "E35" "  -0.95%  "
Set-TextValue "E36" "  -0.26%  "
Set-TextValue "E37" "  +1.13%  "
Set-TextValue "E38" "  -2.65%  "
Set-TextValue "D39" "5.81"
Set-TextValue "E39" "  +0.69%  "
Set-TextValue "D40" "0.987"
Set-TextValue "E40" "  -1.22%  "
Set-TextValue "E41" "  -0.01%  "
Set-TextValue "D43" "47.41"
Set-TextValue "E43" "  -1.62%  "
Set-TextValue "B44" "Arweave"
Set-TextValue "C44" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D44" "43.32"
Set-TextValue "E44" "  +0.12%  "
Set-TextValue "B45" "TheGraph"
Set-TextValue "C45" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D45" "0.298"
Set-TextValue "E45" "  -0.53%  "
Set-TextValue "D46" "152.13"
Set-TextValue "E46" "  +2.93%  "
Set-TextValue "D47" "8.34"
Set-TextValue "E47" "  +0.06%  "
Set-TextValue "D48" "1.86"
Set-TextValue "E48" "  +1.24%  "
Set-TextValue "E49" "  +1.35%  "
Set-TextValue "D50" "389.75"
Set-TextValue "E50" "  -1.22%  "
Set-TextValue "D51" "26.46"
